$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1234603333333333
$ws.Range("H2").Value = 0.370381
$ws.Range("I2").Value = 0.002558470358543426
$ws.Range("J2").Value = 0.002636284444771545
$ws.Range("M2").Value = 26.23985166666667
$ws.Range("N2").Value = 78.719555
$ws.Range("O2").Value = 0.09560625159817936
$ws.Range("P2").Value = 0.09631575414263932
$ws.Range("Q2").Value = 3.239580833383889
$ws.Range("R2").Value = 29.156227500455
$ws.Range("S2").Value = 0.000244605760805387
$ws.Range("T2").Value = 0.0002539157244326806

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1234603333333333
$ws.Range("H3").Value = 0.370381
$ws.Range("I3").Value = 0.002558470358543426
$ws.Range("J3").Value = 0.002636284444771545
$ws.Range("O3").Value = 0.1052038484825964
$ws.Range("P3").Value = 0.1059845756519775
$ws.Range("Q3").Value = 3.564791689301333
$ws.Range("R3").Value = 32.083125203712
$ws.Range("S3").Value = 0.0002691609279474167
$ws.Range("T3").Value = 0.0002794054881770213

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1234603333333333
$ws.Range("H4").Value = 0.370381
$ws.Range("I4").Value = 0.002558470358543426
$ws.Range("J4").Value = 0.002636284444771545
$ws.Range("M4").Value = 108.455335
$ws.Range("N4").Value = 325.366005
$ws.Range("O4").Value = 0.3951626014085634
$ws.Range("P4").Value = 0.3980951384183098
$ws.Range("Q4").Value = 13.38993181087833
$ws.Range("R4").Value = 120.509386297905
$ws.Range("S4").Value = 0.00101101180250872
$ws.Range("T4").Value = 0.001049492020951365

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1234603333333333
$ws.Range("H5").Value = 0.370381
$ws.Range("I5").Value = 0.002558470358543426
$ws.Range("J5").Value = 0.002636284444771545
$ws.Range("M5").Value = 6.0653095
$ws.Range("N5").Value = 12.130619
$ws.Range("O5").Value = 0.02209926768810472
$ws.Range("P5").Value = 0.01484217888683478
$ws.Range("Q5").Value = 0.7488251326398333
$ws.Range("R5").Value = 4.492950795839
$ws.Range("S5").Value = 0.00005654032132553243
$ws.Range("T5").Value = 0.00003912820532587917

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1234603333333333
$ws.Range("H6").Value = 0.370381
$ws.Range("I6").Value = 0.002558470358543426
$ws.Range("J6").Value = 0.002636284444771545
$ws.Range("M6").Value = 104.823008
$ws.Range("N6").Value = 314.469024
$ws.Range("O6").Value = 0.3819280308225561
$ws.Range("P6").Value = 0.3847623529002386
$ws.Range("Q6").Value = 12.94148350868267
$ws.Range("R6").Value = 116.473351578144
$ws.Range("S6").Value = 0.0009771515459563698
$ws.Range("T6").Value = 0.001014343005884599

$ws.Range("I7").Value = 0.9088918061291337
$ws.Range("J7").Value = 0.9365351146153917
$ws.Range("M7").Value = 26.23985166666667
$ws.Range("N7").Value = 78.719555
$ws.Range("O7").Value = 0.09560625159817936
$ws.Range("P7").Value = 0.09631575414263932
$ws.Range("Q7").Value = 1150.855027467237
$ws.Range("R7").Value = 10357.69524720513
$ws.Range("S7").Value = 0.0868957386923056
$ws.Range("T7").Value = 0.09020308584524461

$ws.Range("I8").Value = 0.9088918061291337
$ws.Range("J8").Value = 0.9365351146153917
$ws.Range("O8").Value = 0.1052038484825964
$ws.Range("P8").Value = 0.1059845756519775
$ws.Range("S8").Value = 0.09561891585908278
$ws.Range("T8").Value = 0.0992582767056884

$ws.Range("I9").Value = 0.9088918061291337
$ws.Range("J9").Value = 0.9365351146153917
$ws.Range("M9").Value = 108.455335
$ws.Range("N9").Value = 325.366005
$ws.Range("O9").Value = 0.3951626014085634
$ws.Range("P9").Value = 0.3980951384183098
$ws.Range("Q9").Value = 4756.748213594197
$ws.Range("R9").Value = 42810.73392234778
$ws.Range("S9").Value = 0.3591600505089161
$ws.Range("T9").Value = 0.372830076086422

$ws.Range("I10").Value = 0.9088918061291337
$ws.Range("J10").Value = 0.9365351146153917
$ws.Range("M10").Value = 6.0653095
$ws.Range("N10").Value = 12.130619
$ws.Range("O10").Value = 0.02209926768810472
$ws.Range("P10").Value = 0.01484217888683478
$ws.Range("Q10").Value = 266.0187267783638
$ws.Range("R10").Value = 1596.112360670183
$ws.Range("S10").Value = 0.0200858433231727
$ws.Range("T10").Value = 0.01390022170492396

$ws.Range("I11").Value = 0.9088918061291337
$ws.Range("J11").Value = 0.9365351146153917
$ws.Range("M11").Value = 104.823008
$ws.Range("N11").Value = 314.469024
$ws.Range("O11").Value = 0.3819280308225561
$ws.Range("P11").Value = 0.3847623529002386
$ws.Range("Q11").Value = 4597.437793609419
$ws.Range("R11").Value = 41376.94014248477
$ws.Range("S11").Value = 0.3471312577456565
$ws.Range("T11").Value = 0.3603434542731128

$ws.Range("G12").Value = 4.2730135
$ws.Range("H12").Value = 8.546027
$ws.Range("I12").Value = 0.08854972351232299
$ws.Range("J12").Value = 0.06082860093983664
$ws.Range("M12").Value = 26.23985166666667
$ws.Range("N12").Value = 78.719555
$ws.Range("O12").Value = 0.09560625159817936
$ws.Range("P12").Value = 0.09631575414263932
$ws.Range("Q12").Value = 112.1232404096642
$ws.Range("R12").Value = 672.7394424579851
$ws.Range("S12").Value = 0.00846590714506837
$ws.Range("T12").Value = 0.005858752572962024

$ws.Range("G13").Value = 4.2730135
$ws.Range("H13").Value = 8.546027
$ws.Range("I13").Value = 0.08854972351232299
$ws.Range("J13").Value = 0.06082860093983664
$ws.Range("O13").Value = 0.1052038484825964
$ws.Range("P13").Value = 0.1059845756519775
$ws.Range("Q13").Value = 123.378923430784
$ws.Range("R13").Value = 740.2735405847039
$ws.Range("S13").Value = 0.009315771695566232
$ws.Range("T13").Value = 0.006446893458112065

$ws.Range("G14").Value = 4.2730135
$ws.Range("H14").Value = 8.546027
$ws.Range("I14").Value = 0.08854972351232299
$ws.Range("J14").Value = 0.06082860093983664
$ws.Range("M14").Value = 108.455335
$ws.Range("N14").Value = 325.366005
$ws.Range("O14").Value = 0.3951626014085634
$ws.Range("P14").Value = 0.3980951384183098
$ws.Range("Q14").Value = 463.4311106020225
$ws.Range("R14").Value = 2780.586663612135
$ws.Range("S14").Value = 0.03499153909713858
$ws.Range("T14").Value = 0.02421557031093639

$ws.Range("G15").Value = 4.2730135
$ws.Range("H15").Value = 8.546027
$ws.Range("I15").Value = 0.08854972351232299
$ws.Range("J15").Value = 0.06082860093983664
$ws.Range("M15").Value = 6.0653095
$ws.Range("N15").Value = 12.130619
$ws.Range("O15").Value = 0.02209926768810472
$ws.Range("P15").Value = 0.01484217888683478
$ws.Range("Q15").Value = 25.91714937517825
$ws.Range("R15").Value = 103.668597500713
$ws.Range("S15").Value = 0.001956884043606486
$ws.Range("T15").Value = 0.0009028289765849414

$ws.Range("G16").Value = 4.2730135
$ws.Range("H16").Value = 8.546027
$ws.Range("I16").Value = 0.08854972351232299
$ws.Range("J16").Value = 0.06082860093983664
$ws.Range("M16").Value = 104.823008
$ws.Range("N16").Value = 314.469024
$ws.Range("O16").Value = 0.3819280308225561
$ws.Range("P16").Value = 0.3847623529002386
$ws.Range("Q16").Value = 447.9101282946081
$ws.Range("R16").Value = 2687.460769767648
$ws.Range("S16").Value = 0.03381962153094332
$ws.Range("T16").Value = 0.02340455562124121

